# The "Scénario N" header cells in the recette table are being renumbered:
# Scénario 4 -> 3, 5 -> 4, 6 -> 5, 7 -> 6, 8 -> 7, 9 -> 8, 10 -> 9.
# Word splits the trailing number into its own run (distinct from the
# "Scénario " run) when the author retypes just the digits, so we
# reproduce that: find "Scénario <old>", replace only the trailing
# number characters with the new number, and nudge the run's formatting
# (set then restore the font color) so the engine keeps the freshly
# typed digits as a separate run instead of silently re-merging it back
# into the "Scénario " run it shares identical formatting with.

function Split-ScenarioNumber($d, $oldNumber, $newNumber) {
    $searchText = "Scénario " + $oldNumber
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $numStart = $rng.End - $oldNumber.Length
    $numRng = $d.Range($numStart, $rng.End)
    $numRng.Text = $newNumber
    # force a run split: touch the formatting then put it back
    $numRng.Font.Color = 123456
    $numRng.Font.Color = 16777215
}

$d = $word.ActiveDocument

Split-ScenarioNumber $d "4" "3"
Split-ScenarioNumber $d "5" "4"
Split-ScenarioNumber $d "6" "5"
Split-ScenarioNumber $d "7" "6"
Split-ScenarioNumber $d "8" "7"
# do the already-split "Scénario 9" -> "Scénario 8" before renumbering
# "Scénario 10" -> "Scénario 9", so the Find above still matches a
# unique "Scénario 9" occurrence.
Split-ScenarioNumber $d "9" "8"
Split-ScenarioNumber $d "10" "9"
